# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (strike count) values recalculated for rows 2-22 (column G)
$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 4
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 2
    16 = 2
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
